# Auto-generated edit script applying updated market-price / profit data
# to the per-job Leve tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2182.68
$ws.Range("I98").Value = 2182.68
$ws.Range("K98").Value = 2182.68
$ws.Range("M98").Value = -684.6799999999998

$ws.Range("H112").Value = 2417.875
$ws.Range("J112").Value = 2459.0667
$ws.Range("L112").Value = 7377.2001
$ws.Range("N112").Value = -9593.2001

$ws.Range("H122").Value = 2182.68
$ws.Range("I122").Value = 2182.68
$ws.Range("K122").Value = 6548.039999999999
$ws.Range("M122").Value = -4098.039999999999

$ws.Range("H133").Value = 74150
$ws.Range("J133").Value = 74150
$ws.Range("L133").Value = 74150
$ws.Range("N133").Value = -84270

$ws.Range("H137").Value = 2198.818
$ws.Range("I137").Value = 2103.5217
$ws.Range("K137").Value = 6310.5651
$ws.Range("M137").Value = -3760.5651

$ws.Range("H138").Value = 267255.75
$ws.Range("I138").Value = 5033.05
$ws.Range("J138").Value = 337181.8
$ws.Range("K138").Value = 15099.15
$ws.Range("L138").Value = 1011545.4
$ws.Range("M138").Value = -9959.150000000001
$ws.Range("N138").Value = -1021825.4

$ws.Range("H141").Value = 1500
$ws.Range("I141").Value = 1500
$ws.Range("K141").Value = 4500
$ws.Range("M141").Value = 680

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2153.14
$ws.Range("I32").Value = 2153.14
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2153.14
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1866.14
$ws.Range("N32").ClearContents()

$ws.Range("H61").Value = 6503.778
$ws.Range("I61").Value = 5642.5454
$ws.Range("J61").Value = 7857.143
$ws.Range("K61").Value = 5642.5454
$ws.Range("L61").Value = 7857.143
$ws.Range("M61").Value = -5430.5454
$ws.Range("N61").Value = -8281.143

$ws.Range("H74").Value = 560094.2
$ws.Range("I74").Value = 698930.25
$ws.Range("K74").Value = 698930.25
$ws.Range("M74").Value = -698056.25

$ws.Range("H77").Value = 560094.2
$ws.Range("I77").Value = 698930.25
$ws.Range("K77").Value = 3494651.25
$ws.Range("M77").Value = -3490283.25

$ws.Range("H102").Value = 4365.1816
$ws.Range("I102").Value = 4465.067
$ws.Range("K102").Value = 4465.067
$ws.Range("M102").Value = -2843.067

$ws.Range("H122").Value = 3327.3389
$ws.Range("I122").Value = 3067.9058
$ws.Range("J122").Value = 5619
$ws.Range("K122").Value = 9203.7174
$ws.Range("L122").Value = 16857
$ws.Range("M122").Value = -6753.7174
$ws.Range("N122").Value = -21757

$ws.Range("H132").Value = 3033.814
$ws.Range("I132").Value = 2686.2812
$ws.Range("J132").Value = 4044.818
$ws.Range("K132").Value = 8058.8436
$ws.Range("L132").Value = 12134.454
$ws.Range("M132").Value = -5528.8436
$ws.Range("N132").Value = -17194.454

$ws.Range("H136").Value = 6503.778
$ws.Range("I136").Value = 5642.5454
$ws.Range("J136").Value = 7857.143
$ws.Range("K136").Value = 16927.6362
$ws.Range("L136").Value = 23571.429
$ws.Range("M136").Value = -14377.6362
$ws.Range("N136").Value = -28671.429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4752.467
$ws.Range("I99").Value = 3473.8572
$ws.Range("K99").Value = 3473.8572
$ws.Range("M99").Value = -1975.8572

$ws.Range("H105").Value = 23638692
$ws.Range("I105").Value = 1430044.6
$ws.Range("K105").Value = 1430044.6
$ws.Range("M105").Value = -1428297.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4025.9
$ws.Range("I31").Value = 3381.9124
$ws.Range("J31").Value = 5138.242
$ws.Range("K31").Value = 3381.9124
$ws.Range("L31").Value = 5138.242
$ws.Range("M31").Value = -3086.9124
$ws.Range("N31").Value = -5728.242

$ws.Range("H34").Value = 4025.9
$ws.Range("I34").Value = 3381.9124
$ws.Range("J34").Value = 5138.242
$ws.Range("K34").Value = 3381.9124
$ws.Range("L34").Value = 5138.242
$ws.Range("M34").Value = -3179.9124
$ws.Range("N34").Value = -5542.242

$ws.Range("H58").Value = 3064.1538
$ws.Range("I58").Value = 2231.25
$ws.Range("J58").Value = 3778.0715
$ws.Range("K58").Value = 2231.25
$ws.Range("L58").Value = 3778.0715
$ws.Range("M58").Value = -2028.25
$ws.Range("N58").Value = -4184.0715

$ws.Range("H132").Value = 3565.7297
$ws.Range("I132").Value = 3380.2
$ws.Range("K132").Value = 10140.6
$ws.Range("M132").Value = -7610.599999999999

$ws.Range("H136").Value = 3064.1538
$ws.Range("I136").Value = 2231.25
$ws.Range("J136").Value = 3778.0715
$ws.Range("K136").Value = 6693.75
$ws.Range("L136").Value = 11334.2145
$ws.Range("M136").Value = -4143.75
$ws.Range("N136").Value = -16434.2145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1993.3704
$ws.Range("I5").Value = 431.17648
$ws.Range("J5").Value = 4649.1
$ws.Range("K5").Value = 1293.52944
$ws.Range("L5").Value = 13947.3
$ws.Range("M5").Value = -1181.52944
$ws.Range("N5").Value = -14171.3

$ws.Range("H86").Value = 268.57144
$ws.Range("I86").Value = 268.57144
$ws.Range("K86").Value = 805.71432
$ws.Range("M86").Value = 380.28568

$ws.Range("H89").Value = 268.57144
$ws.Range("I89").Value = 268.57144
$ws.Range("K89").Value = 2417.14296
$ws.Range("M89").Value = 3510.85704

$ws.Range("H97").Value = 650999
$ws.Range("I97").Value = 1666998.4
$ws.Range("J97").Value = 41399.4
$ws.Range("K97").Value = 5000995.199999999
$ws.Range("L97").Value = 124198.2
$ws.Range("M97").Value = -5000499.199999999
$ws.Range("N97").Value = -125190.2

$ws.Range("H113").Value = 2054.1667
$ws.Range("J113").Value = 2431.7778
$ws.Range("L113").Value = 7295.3334
$ws.Range("N113").Value = -11635.3334

$ws.Range("H120").Value = 30000
$ws.Range("J120").Value = 30000
$ws.Range("L120").Value = 90000
$ws.Range("N120").Value = -99676

$ws.Range("H122").Value = 1529.8889
$ws.Range("J122").Value = 1541.2858
$ws.Range("L122").Value = 13871.5722
$ws.Range("N122").Value = -18771.5722

$ws.Range("H131").Value = 5631.927
$ws.Range("J131").Value = 2181.0264
$ws.Range("L131").Value = 6543.0792
$ws.Range("N131").Value = -16623.0792

$ws.Range("H132").Value = 3560.12
$ws.Range("I132").Value = 2641.8
$ws.Range("J132").Value = 3789.7
$ws.Range("K132").Value = 23776.2
$ws.Range("L132").Value = 34107.3
$ws.Range("M132").Value = -21246.2
$ws.Range("N132").Value = -39167.3

$ws.Range("H135").Value = 1993.3704
$ws.Range("I135").Value = 431.17648
$ws.Range("J135").Value = 4649.1
$ws.Range("K135").Value = 3880.58832
$ws.Range("L135").Value = 41841.9
$ws.Range("M135").Value = -1345.58832
$ws.Range("N135").Value = -46911.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2066.4075
$ws.Range("I102").Value = 1922.7778
$ws.Range("K102").Value = 1922.7778
$ws.Range("M102").Value = -300.7778000000001

$ws.Range("H113").Value = 6279.5454
$ws.Range("I113").Value = 6407.6
$ws.Range("K113").Value = 6407.6
$ws.Range("M113").Value = -4237.6

$ws.Range("H122").Value = 2511.2727
$ws.Range("I122").Value = 2561.6
$ws.Range("J122").Value = 2008
$ws.Range("K122").Value = 7684.799999999999
$ws.Range("L122").Value = 6024
$ws.Range("M122").Value = -5234.799999999999
$ws.Range("N122").Value = -10924

$ws.Range("H132").Value = 3291.4614
$ws.Range("I132").Value = 2897.4062
$ws.Range("J132").Value = 5092.857
$ws.Range("K132").Value = 8692.2186
$ws.Range("L132").Value = 15278.571
$ws.Range("M132").Value = -6162.2186
$ws.Range("N132").Value = -20338.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6179.885
$ws.Range("I7").Value = 4550.7334
$ws.Range("K7").Value = 4550.7334
$ws.Range("M7").Value = -4438.7334

$ws.Range("H22").Value = 119050030
$ws.Range("I22").Value = 14287343
$ws.Range("K22").Value = 14287343
$ws.Range("M22").Value = -14287048

$ws.Range("H27").Value = 119050030
$ws.Range("I27").Value = 14287343
$ws.Range("K27").Value = 14287343
$ws.Range("M27").Value = -14287236

$ws.Range("H46").Value = 1491.7441
$ws.Range("I46").Value = 2498.75
$ws.Range("K46").Value = 2498.75
$ws.Range("M46").Value = -2310.75

$ws.Range("H68").Value = 3080.6365
$ws.Range("I68").Value = 2757
$ws.Range("J68").Value = 3647
$ws.Range("K68").Value = 2757
$ws.Range("L68").Value = 3647
$ws.Range("M68").Value = -2008
$ws.Range("N68").Value = -5145

$ws.Range("H71").Value = 3080.6365
$ws.Range("I71").Value = 2757
$ws.Range("J71").Value = 3647
$ws.Range("K71").Value = 13785
$ws.Range("L71").Value = 18235
$ws.Range("M71").Value = -10041
$ws.Range("N71").Value = -25723

$ws.Range("H122").Value = 7875.5
$ws.Range("I122").Value = 7857
$ws.Range("K122").Value = 23571
$ws.Range("M122").Value = -21121

$ws.Range("H126").Value = 6179.885
$ws.Range("I126").Value = 4550.7334
$ws.Range("K126").Value = 13652.2002
$ws.Range("M126").Value = -11182.2002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4339.684
$ws.Range("I81").Value = 3479.4
$ws.Range("K81").Value = 6958.8
$ws.Range("M81").Value = -5897.8

$ws.Range("H84").Value = 4339.684
$ws.Range("I84").Value = 3479.4
$ws.Range("K84").Value = 34794
$ws.Range("M84").Value = -29490

$ws.Range("H122").Value = 11907941
$ws.Range("I122").Value = 2716.7334
$ws.Range("K122").Value = 8150.2002
$ws.Range("M122").Value = -5700.2002
